$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date header for column AR (next day after AQ1 = "07-ago")
$ws.Range("AR1").Value = "08-ago"
$ws.Range("AR1").NumberFormat = "@"

# New data values for column AR, rows 2-11 (matching the numeric, centered
# style already used by the rest of the date columns)
$values = @(14, 14, 11, 16, 8, 16, 17, 20, 12, 21)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("AR$row")
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
}

# Move the active selection to match the post-edit state
$ws.Range("AR12").Select()
